$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (TSM)
$ws.Range("D2").Value = 292.09
$ws.Range("F2").Value = 2.62
$ws.Range("K2").Value = 63
$ws.Range("N2").Value = 66.04328690552585

# Row 3 (ASML)
$ws.Range("D3").Value = 1108.78
$ws.Range("E3").Value = 63.3
$ws.Range("F3").Value = 12.25
$ws.Range("J3").Value = 63
$ws.Range("K3").Value = 56
$ws.Range("N3").Value = 66.04328690552585

# Row 4 (QCOM)
$ws.Range("D4").Value = 170.7
$ws.Range("E4").Value = 45.5
$ws.Range("F4").Value = 3.42
$ws.Range("K4").Value = 55.8
$ws.Range("N4").Value = 66.04328690552585

# Row 5 (NVDA)
$ws.Range("D5").Value = 181.46
$ws.Range("F5").Value = -0.6
$ws.Range("K5").Value = 49.2
$ws.Range("N5").Value = 66.04328690552585

# Row 6 (AMD)
$ws.Range("D6").Value = 215.24
$ws.Range("E6").Value = 40.3
$ws.Range("F6").Value = 0.09
$ws.Range("J6").Value = 60
$ws.Range("K6").Value = 44.8
$ws.Range("N6").Value = 66.04328690552585
